$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on cells whose new values look numeric,
# so Excel keeps the original formatted text (e.g. trailing zeros)
# instead of auto-converting to a floating point number.
$ws.Range("D2").Value = '28.561.36'
$ws.Range("E2").Value = '  +0.17%  '
$ws.Range("D3").Value = '1.870.88'
$ws.Range("E3").Value = '  -0.79%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.010'
$ws.Range("E4").Value = '  -0.82%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.84'
$ws.Range("E5").Value = '  -0.32%  '
$ws.Range("E6").Value = '  -0.57%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5071'
$ws.Range("E7").Value = '  -1.72%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3892'
$ws.Range("E8").Value = '  -1.57%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08365'
$ws.Range("E9").Value = '  +0.27%  '
$ws.Range("B10").Value = 'Polygon'
$ws.Range("C10").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.103'
$ws.Range("E10").Value = '  -2.07%  '
$ws.Range("B11").Value = 'OKB'
$ws.Range("C11").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '41.80'
$ws.Range("E11").Value = '  -0.57%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.201'
$ws.Range("E12").Value = '  -1.65%  '
$ws.Range("D13").Value = '1.873.60'
$ws.Range("E13").Value = '  +0.43%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.39'
$ws.Range("E14").Value = '  -0.45%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.222'
$ws.Range("E15").Value = '  -0.79%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.010'
$ws.Range("E16").Value = '  -0.81%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001099'
$ws.Range("E17").Value = '  -0.88%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '90.95'
$ws.Range("E18").Value = '  -0.67%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06709'
$ws.Range("E19").Value = '  -0.64%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.68'
$ws.Range("E20").Value = '  -0.76%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.009'
$ws.Range("E21").Value = '  -0.53%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.914'
$ws.Range("E22").Value = '  -1.89%  '
$ws.Range("D23").Value = '28.577.24'
$ws.Range("E23").Value = '  +0.26%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.04'
$ws.Range("E24").Value = '  -1.42%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.238'
$ws.Range("E25").Value = '  -1.11%  '
$ws.Range("D26").Value = '2.086.57'
$ws.Range("E26").Value = '  +0.26%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '161.59'
$ws.Range("E27").Value = '  +0.11%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.63'
$ws.Range("E28").Value = '  -1.22%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.342'
$ws.Range("E29").Value = '  -4.72%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '125.38'
$ws.Range("E30").Value = '  -1.49%  '
$ws.Range("E31").Value = '  -2.42%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.040'
$ws.Range("E32").Value = '  -1.59%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.783'
$ws.Range("E33").Value = '  -2.70%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.610'
$ws.Range("E34").Value = '  -1.16%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02446'
$ws.Range("E35").Value = '  -0.53%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06521'
$ws.Range("E36").Value = '  -0.81%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2158'
$ws.Range("E37").Value = '  -1.88%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.859'
$ws.Range("E38").Value = '  -4.58%  '
$ws.Range("B39").Value = 'InternetComputer(DFINITY)'
$ws.Range("C39").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.040'
$ws.Range("E39").Value = '  +0.54%  '
$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.254'
$ws.Range("E40").Value = '  -0.40%  '
$ws.Range("E41").Value = '  -0.60%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6415'
$ws.Range("E42").Value = '  -1.45%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '11.08'
$ws.Range("E43").Value = '  -1.40%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.009'
$ws.Range("E44").Value = '  -0.44%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6000'
$ws.Range("E45").Value = '  -1.61%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '13.00'
$ws.Range("E46").Value = '  -1.01%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.687'
$ws.Range("E47").Value = '  -0.87%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.996'
$ws.Range("E48").Value = '  -1.81%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.211'
$ws.Range("E49").Value = '  -1.26%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '121.67'
$ws.Range("E50").Value = '  -0.76%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.181'
$ws.Range("E51").Value = '  -8.31%  '
